$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G represents strikeout count "K" for each start (row 1 is header).
# Regenerated values replace the previous (incorrect "Strike#") figures.
$kValues = @{
    2 = 2
    3 = 6
    4 = 3
    5 = 2
    6 = 3
    7 = 2
    8 = 4
    9 = 4
    10 = 4
    11 = 6
    12 = 4
    13 = 7
    14 = 7
    15 = 4
    16 = 3
    17 = 5
    18 = 3
    19 = 7
    20 = 6
    21 = 5
    22 = 3
    23 = 8
    24 = 5
    25 = 5
    26 = 12
    27 = 8
    28 = 6
    29 = 5
    30 = 6
    31 = 6
    32 = 7
    33 = 4
    34 = 7
    35 = 1
    36 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
